$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 931
$ws.Range("I38").Value = 86.416664
$ws.Range("J38").Value = 1944.5
$ws.Range("K38").Value = 259.249992
$ws.Range("L38").Value = 5833.5
$ws.Range("M38").Value = 112.750008
$ws.Range("N38").Value = -6577.5

$ws.Range("H39").Value = 594.63635
$ws.Range("I39").Value = 595
$ws.Range("K39").Value = 1785
$ws.Range("M39").Value = -1489

$ws.Range("H129").Value = 894.1818
$ws.Range("J129").Value = 899.62964
$ws.Range("L129").Value = 2698.88892
$ws.Range("N129").Value = -12698.88892

$ws.Range("H132").Value = 8551166
$ws.Range("I132").Value = 10419030
$ws.Range("J132").Value = 12357.286
$ws.Range("K132").Value = 31257090
$ws.Range("L132").Value = 37071.858
$ws.Range("M132").Value = -31254560
$ws.Range("N132").Value = -42131.858

$ws.Range("H137").Value = 1160.6578
$ws.Range("J137").Value = 1531.9166
$ws.Range("L137").Value = 4595.7498
$ws.Range("N137").Value = -9695.7498

$ws.Range("H138").Value = 2089.6736
$ws.Range("I138").Value = 1075.091
$ws.Range("J138").Value = 2217.954
$ws.Range("K138").Value = 3225.273
$ws.Range("L138").Value = 6653.862000000001
$ws.Range("M138").Value = 1914.727
$ws.Range("N138").Value = -16933.862

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2049.0488
$ws.Range("I32").Value = 2200.081
$ws.Range("J32").Value = 652
$ws.Range("K32").Value = 2200.081
$ws.Range("L32").Value = 652
$ws.Range("M32").Value = -1913.081
$ws.Range("N32").Value = -1226

$ws.Range("H45").Value = 2318.3333
$ws.Range("I45").Value = 1930
$ws.Range("J45").Value = 2928.5715
$ws.Range("K45").Value = 1930
$ws.Range("L45").Value = 2928.5715
$ws.Range("M45").Value = -1553
$ws.Range("N45").Value = -3682.5715

$ws.Range("H55").Value = 37800
$ws.Range("J55").Value = 37800
$ws.Range("L55").Value = 37800
$ws.Range("N55").Value = -38430

$ws.Range("H61").Value = 1422.5454
$ws.Range("I61").Value = 1250.9474
$ws.Range("K61").Value = 1250.9474
$ws.Range("M61").Value = -1038.9474

$ws.Range("H118").Value = 38000
$ws.Range("J118").Value = 38000
$ws.Range("L118").Value = 38000
$ws.Range("N118").Value = -41314

$ws.Range("H132").Value = 3060.25
$ws.Range("J132").Value = 2912.4443
$ws.Range("L132").Value = 8737.332900000001
$ws.Range("N132").Value = -13797.3329

$ws.Range("H136").Value = 1422.5454
$ws.Range("I136").Value = 1250.9474
$ws.Range("K136").Value = 3752.8422
$ws.Range("M136").Value = -1202.8422

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 19000
$ws.Range("I26").Value = 19000
$ws.Range("K26").Value = 19000
$ws.Range("M26").Value = -18708

$ws.Range("H116").Value = 28944.5
$ws.Range("J116").Value = 28944.5
$ws.Range("L116").Value = 28944.5
$ws.Range("N116").Value = -38122.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 83334584
$ws.Range("I16").Value = 142858300
$ws.Range("J16").Value = 1378.6
$ws.Range("K16").Value = 142858300
$ws.Range("L16").Value = 1378.6
$ws.Range("M16").Value = -142858013
$ws.Range("N16").Value = -1952.6

$ws.Range("H41").Value = 9546
$ws.Range("J41").Value = 23500
$ws.Range("L41").Value = 23500
$ws.Range("N41").Value = -24356

$ws.Range("H113").Value = 83334584
$ws.Range("I113").Value = 142858300
$ws.Range("J113").Value = 1378.6
$ws.Range("K113").Value = 142858300
$ws.Range("L113").Value = 1378.6
$ws.Range("M113").Value = -142856130
$ws.Range("N113").Value = -5718.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1472.625
$ws.Range("J68").Value = 2303.3333
$ws.Range("L68").Value = 6909.999899999999
$ws.Range("N68").Value = -8531.999899999999

$ws.Range("H71").Value = 1472.625
$ws.Range("J71").Value = 2303.3333
$ws.Range("L71").Value = 20729.9997
$ws.Range("N71").Value = -28841.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 19095
$ws.Range("J26").Value = 19095
$ws.Range("L26").Value = 19095
$ws.Range("N26").Value = -19655

$ws.Range("H45").Value = 37541.668
$ws.Range("J45").Value = 37541.668
$ws.Range("L45").Value = 37541.668
$ws.Range("N45").Value = -38659.668

$ws.Range("H50").Value = 19095
$ws.Range("J50").Value = 19095
$ws.Range("L50").Value = 19095
$ws.Range("N50").Value = -20091

$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").Value = $null

$ws.Range("H107").Value = 716.8823
$ws.Range("I107").Value = 593.36365
$ws.Range("J107").Value = 943.3333
$ws.Range("K107").Value = 593.36365
$ws.Range("L107").Value = 943.3333
$ws.Range("M107").Value = 1326.63635
$ws.Range("N107").Value = -4783.3333

$ws.Range("H113").Value = 2704.8333
$ws.Range("I113").Value = 1245.8
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 1245.8
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = 924.2
$ws.Range("N113").Value = -14340

$ws.Range("H117").Value = 28000
$ws.Range("J117").Value = 28000
$ws.Range("L117").Value = 28000
$ws.Range("N117").Value = -34884

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = $null
$ws.Range("N61").Value = $null

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 0
$ws.Range("N113").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 666
$ws.Range("I113").Value = 332.5
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 997.5
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = 1172.5
$ws.Range("N113").Value = -10340

$ws.Range("H122").Value = 7880300
$ws.Range("I122").Value = 8966952
$ws.Range("K122").Value = 26900856
$ws.Range("M122").Value = -26898406
